$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bigval = 1.608817645280888 * [Math]::Pow(10, 37)

# Row 2
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 19.48425592650926

# Row 3
$ws.Range("B3").Value = 0.6545652718822623
$ws.Range("C3").Value = 22898927661.19195
$ws.Range("D3").Value = 186123.597850132
$ws.Range("E3").Value = $bigval
$ws.Range("G3").Value = $bigval

# Row 4
$ws.Range("B4").Value = 0.6545652718822623
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 13.86384647080068
$ws.Range("G4").Value = 16.86649396021207

# Row 5
$ws.Range("B5").Value = 1.445647641019636
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 13.86384647080068
$ws.Range("G5").Value = 17.65757632934944

# Row 6
$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 0.7210945179870265
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 6.15379541431027

# Row 7
$ws.Range("B7").Value = 0.2881169905109251
$ws.Range("C7").Value = 0.3048912486333797
$ws.Range("D7").Value = 3.223369029078222
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("G7").Value = 4.349763226824225
